# Update the "ランサーズ" (Lancers) worksheet with the latest scrape run
# (2025-12-05 06:28:47 JST), replacing the previous run's rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Drop the old rows 8-23 (only 6 fresh listings remain this run) ---
$ws.Rows("8:23").Delete()

# --- 2. Overwrite rows 2-7 with the new scrape data ---
$timestamp = "2025-12-05 06:28:47"

$ws.Range("A2").Value = $timestamp
$ws.Range("B2").Value = "webアプリの開発"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5437832"
$ws.Range("G2").Value = 100
$ws.Range("H2").Value = "◆開発 ◇アプリ"

$ws.Range("A3").Value = $timestamp
$ws.Range("B3").Value = "Dify、RAGシステムの相談が可能な方を募集!疑問や課題を解決してくれる相談役【リモート作業OK】"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5447772"
$ws.Range("G3").Value = 28
$ws.Range("H3").ClearContents()

$ws.Range("A4").Value = $timestamp
$ws.Range("B4").Value = "Accessシステムのバージョンup対応(Access2010→Access2021)"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "20,000 円 ~ 30,000 円 / 募集期間 2 日、取引期間 0 日"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5447970"
$ws.Range("G4").Value = 25
$ws.Range("H4").ClearContents()

$ws.Range("A5").Value = $timestamp
$ws.Range("B5").Value = "〖リモート可〗Delphiエンジニア募集"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5341051"
$ws.Range("G5").Value = 25
$ws.Range("H5").ClearContents()

$ws.Range("A6").Value = $timestamp
$ws.Range("B6").Value = "Azure SQL Databaseの登録内容をHPに表示(絞り込み検索・フォームからDB登録)"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5447771"
$ws.Range("G6").Value = 18
$ws.Range("H6").ClearContents()

$ws.Range("A7").Value = $timestamp
$ws.Range("B7").Value = "【急募】ドメインメールの設定サポートをお願いします"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5447991"
$ws.Range("G7").Value = 10
$ws.Range("H7").ClearContents()

# --- 3. Rebuild the hyperlinks so only F2:F7 point at the new URLs ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5437832")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5447772")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5447970")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5341051")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5447771")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5447991")

# Re-apply the "Hyperlink" cell style so F2:F7 reuse the workbook's existing
# Hyperlink style slot instead of a freshly duplicated one.
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("F5").Style = "Hyperlink"
$ws.Range("F6").Style = "Hyperlink"
$ws.Range("F7").Style = "Hyperlink"

# --- 4. Column width tweaks (B, D, H). COM ColumnWidth reports ~0.83
#        narrower than the stored XML width, so compensate by that offset
#        to land on the exact target widths of 52 / 41 / 12. ---
$ws.Columns("B").ColumnWidth = 52 - 5/6
$ws.Columns("D").ColumnWidth = 41 - 5/6
$ws.Columns("H").ColumnWidth = 12 - 5/6
